$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Anzahl aller moeglichen Loesungen" column (C) with new
# algorithm results for Aufgabe 3.
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 8
$ws.Range("C8").Value = 1040

# N = 6 no longer has a total-solutions figure listed.
$ws.Range("C9").ClearContents()

# Update the "mit einer Reihe [...]" column (G) counts.
$ws.Range("G8").Value = 28
$ws.Range("G9").Value = 12

# Move the selection off the old header range onto the next empty row.
$ws.Range("E17").Select()
